# The presentation's slide-master theme (theme1.xml, "Integral" / "Red Violet")
# is switched over to the standard default "Office Theme" color scheme
# ("Office" clrScheme): dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink all change to
# the stock Office palette. Fonts and the format scheme are already identical
# between the two themes, so only the twelve theme colors need to move.
#
# PowerPoint's ThemeColorScheme collection addresses the slots in the fixed
# OOXML clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# .RGB takes/returns a packed 0xBBGGRR integer (VBA RGB() order), so the
# values below are the target srgbClr hex bytes reversed.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$cs = $master.Theme.ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  -> 4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  -> 70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    -> 0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink -> 954F72
